$wb = $excel.ActiveWorkbook

# "zh-cn" worksheet: row 16 corresponds to file
# add3c42c-3d0a-49af-9bb5-416741cad363.88cf6a626098b018e0cf9590eeef0a43d1d1ee7b.zh-cn.xlf
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D16").Value = "2016-03-08 02:37:54"
$wsZhCn.Range("G16").Value = "2016-03-08 02:38:37"

# "de-de" worksheet: row 16 corresponds to file
# add3c42c-3d0a-49af-9bb5-416741cad363.88cf6a626098b018e0cf9590eeef0a43d1d1ee7b.de-de.xlf
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D16").Value = "2016-03-08 02:38:02"
$wsDeDe.Range("G16").Value = "2016-03-08 02:38:51"
